$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row insertions -------------------------------------------------
# 1) Two rows inserted at row 7: the first becomes the new "ERR104 / Send
#    message failed..." entry, the second is a blank spacer row (matching
#    the blank-row-per-block layout used throughout the sheet). Everything
#    that used to be at row 7+ shifts down by 2.
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# 2) One more row inserted at (the now-shifted) row 55: this becomes the
#    new "WARN702 / Shown if File>download report..." entry. The blank
#    spacer row that used to sit at (shifted) row 55 is pushed to row 56,
#    and everything below shifts down by 1 more (net +3 vs. the original).
$ws.Rows.Item(55).Insert()

# --- New cell content -------------------------------------------------
# Populate the second insertion point first so the new shared strings are
# appended to sharedStrings.xml in the same order as the target file
# (index 154, 155, 156, 157).
$ws.Range("A55").Value = "Shown if File>download report is pressed and no data collection related to the user was found"
$ws.Range("B55").Value = "WARN702"
$ws.Range("C55").Value = "yes"

$ws.Range("B7").Value = "ERR104"
$ws.Range("A7").Value = "Send message failed as result of SOAP call (probably due to username wrongly typed)"
$ws.Range("C7").Value = "yes"

# --- Styling ------------------------------------------------------------
# The new row 55 picks up the same border styling used by the other
# bordered A/B pairs in this sheet (e.g. row 53, row 66) - style index 1
# for column A, 2 for column B. Copy the formatting only (not the value).
$ws.Range("A53").Copy()
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("B53").Copy()
$ws.Range("B55").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Table range ----------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C66"))

# --- View state -------------------------------------------------------
$ws.Range("A8").Select()
